$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13 (columns A:D)
$data = @(
    @(2, 1, 5, 6),
    @(5, 1, 11, 11),
    @(6, 1, 16, 16),
    @(7, 1, 21, 21),
    @(8, 1, 26, 26),
    @(9, 1, 31, 31),
    @(9, 1, 36, 36),
    @(10, 2, 5, 6),
    @(1, 3, 5, 5),
    @(4, 3, 10, 10),
    @(3, 4, 5, 5),
    @(11, 4, 10, 10)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
